$d = $word.ActiveDocument

$replacements = @(
    @("87÷8=", "47÷4="),
    @("63÷8=", "31÷7="),
    @("16÷7=", "77÷4="),
    @("97÷6=", "89÷5="),
    @("51÷3=", "13÷4="),
    @("64÷2=", "91÷8="),
    @("32÷5=", "35÷7="),
    @("97÷9=", "75÷2="),
    @("45÷2=", "72÷2="),
    @("46÷5=", "63÷2="),
    @("43÷4=", "71÷4="),
    @("90÷8=", "43÷7="),
    @("16÷6=", "56÷9="),
    @("19÷9=", "66÷7="),
    @("79÷2=", "86÷2="),
    @("89÷4=", "40÷2="),
    @("49÷6=", "40÷7="),
    @("22÷8=", "35÷2="),
    @("11÷6=", "43÷9="),
    @("24÷7=", "78÷3="),
    @("65÷7=", "42÷4="),
    @("53÷9=", "73÷7="),
    @("49÷8=", "69÷4="),
    @("13÷8=", "53÷7="),
    @("73÷2=", "60÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
